# Apply the HackerRank progress-tracker update:
#  - Fill in row 13 (3rd Problem Solving entry) with new data
#  - Stamp the "Problem Solving" category label down through rows 14-17
#  - Move the active selection to N16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: fill in the newly logged Problem-Solving submission ---
$ws.Range("M13").Value = "Problem Solving(Algorithms & Data Structures)"
$ws.Range("N13").Value = 45124
$ws.Range("O13").Value = "1207.97/2200"
$ws.Range("P13").Value = 118419
$ws.Range("Q13").Formula = "=IF(ROW()>2,(`$P`$2-P13)/`$P`$2,""NA"")"

# --- Rows 14-17: carry the category label down (values/formats only; the
#     rest of each row is still empty, waiting to be filled later) ---
$ws.Range("M14").Value = "Problem Solving(Algorithms & Data Structures)"
$ws.Range("M15").Value = "Problem Solving(Algorithms & Data Structures)"
$ws.Range("M16").Value = "Problem Solving(Algorithms & Data Structures)"

# Row 17 previously used a different style (no fill/border matching the
# table body); copy the formatting from M13 before setting its value so it
# matches the rest of the category column.
$ws.Range("M13").Copy()
$ws.Range("M17").PasteSpecial(-4122)
$ws.Range("M17").Value = "Problem Solving(Algorithms & Data Structures)"
$excel.CutCopyMode = 0

# --- Move the selection to reflect where the user ended up editing ---
$ws.Range("N16").Select()
